$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.947.13'
$ws.Range('E2').Value = '  +0.56%  '
$ws.Range('D3').Value = '2.301.56'
$ws.Range('E3').Value = '  -0.05%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = '''309.84'
$ws.Range('E5').Value = '  -2.19%  '
$ws.Range('D6').Value = '''104.95'
$ws.Range('E6').Value = '  +1.09%  '
$ws.Range('D7').Value = '''0.625'
$ws.Range('E7').Value = '  -0.51%  '
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('D9').Value = '''0.604'
$ws.Range('E9').Value = '  -0.09%  '
$ws.Range('D10').Value = '''39.59'
$ws.Range('E10').Value = '  -0.58%  '
$ws.Range('E11').Value = '  +0.31%  '
$ws.Range('E12').Value = '  -2.65%  '
$ws.Range('E13').Value = '  -0.06%  '
$ws.Range('D14').Value = '''0.983'
$ws.Range('E14').Value = '  -1.61%  '
$ws.Range('D15').Value = '''15.25'
$ws.Range('E15').Value = '  -1.01%  '
$ws.Range('D16').Value = '2.651.80'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '2.313.25'
$ws.Range('E17').Value = '  +0.54%  '
$ws.Range('D18').Value = '42.811.41'
$ws.Range('E18').Value = '  +0.43%  '
$ws.Range('D19').Value = '''7.34'
$ws.Range('E19').Value = '  -3.16%  '
$ws.Range('D20').Value = '''13.79'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  -1.23%  '
$ws.Range('D22').Value = '''73.36'
$ws.Range('E22').Value = '  -0.95%  '
$ws.Range('E23').Value = '  -2.74%  '
$ws.Range('D24').Value = '''267.87'
$ws.Range('E24').Value = '  +0.53%  '
$ws.Range('D25').Value = '''2.25'
$ws.Range('E25').Value = '  +0.27%  '
$ws.Range('D26').Value = '''1.01'
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('D27').Value = '''7.79'
$ws.Range('E27').Value = '  +17.42%  '
$ws.Range('D28').Value = '''10.93'
$ws.Range('E28').Value = '  +0.14%  '
$ws.Range('E29').Value = '  -2.02%  '
$ws.Range('D30').Value = '''37.77'
$ws.Range('E30').Value = '  +0.27%  '
$ws.Range('D31').Value = '''22.17'
$ws.Range('E31').Value = '  -2.08%  '
$ws.Range('D32').Value = '''165.97'
$ws.Range('E32').Value = '  +0.10%  '
$ws.Range('D33').Value = '''0.0865'
$ws.Range('E33').Value = '  -2.43%  '
$ws.Range('D34').Value = '''2.82'
$ws.Range('E34').Value = '  +5.62%  '
$ws.Range('E35').Value = '  -0.93%  '
$ws.Range('E36').Value = '  -0.68%  '
$ws.Range('D37').Value = '''4.62'
$ws.Range('E37').Value = '  +0.90%  '
$ws.Range('E38').Value = '  +1.23%  '
$ws.Range('E39').Value = '  +1.05%  '
$ws.Range('D40').Value = '''3.60'
$ws.Range('E40').Value = '  -3.91%  '
$ws.Range('D41').Value = '''107.97'
$ws.Range('E41').Value = '  +13.37%  '
$ws.Range('E42').Value = '  -3.79%  '
$ws.Range('D43').Value = '''71.12'
$ws.Range('E43').Value = '  +0.64%  '
$ws.Range('E44').Value = '  +0.46%  '
$ws.Range('D45').Value = '''1.01'
$ws.Range('E45').Value = '  +0.30%  '
$ws.Range('D46').Value = '''12.23'
$ws.Range('E46').Value = '  -2.10%  '
$ws.Range('D47').Value = '1.696.93'
$ws.Range('E47').Value = '  +2.35%  '
$ws.Range('D48').Value = '''111.29'
$ws.Range('E48').Value = '  -5.67%  '
$ws.Range('D49').Value = '''75.84'
$ws.Range('D50').Value = '''8.85'
$ws.Range('E50').Value = '  -0.64%  '
$ws.Range('D51').Value = '''5.16'
$ws.Range('E51').Value = '  -2.51%  '
